$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the comparative period year typo: 2011 -> 2021
$ws.Range("C2").Value = "Dec 1, 2021 - Dec 1, 2026"

# Move the active selection to C2, matching where the user last edited
$ws.Range("C2").Select()
